# Weekly update to fruit/vegetable price data: swap the data between
# row 2 <-> row 4, and row 3 <-> row 5 (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($row1, $row2, $col) {
    $r1 = $ws.Range($col + $row1)
    $r2 = $ws.Range($col + $row2)
    $v1 = $r1.Value()
    $v2 = $r2.Value()
    $r1.Value = $v2
    $r2.Value = $v1
}

$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    Swap-Cell 2 4 $col
}

foreach ($col in $cols) {
    Swap-Cell 3 5 $col
}
